# Apply updated "dSF" (column F) values as per the commit:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = 5
    "F4"  = -1
    "F7"  = -3
    "F8"  = -1
    "F9"  = 1
    "F10" = -4
    "F11" = 4
    "F13" = 7
    "F14" = 2
    "F17" = 3
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
